$wb = $excel.ActiveWorkbook

# ---- Sheet: RUNMANAGER ----
$ws1 = $wb.Worksheets.Item("RUNMANAGER")
$ws1.Range("C4").Value = "yes"
$ws1.Range("C7").Value = "no"
[void]$ws1.Range("C4").Select()

# ---- Sheet: DATA ----
$ws2 = $wb.Worksheets.Item("DATA")

# Insert a new "version" column after the "browser" column (D)
[void]$ws2.Columns.Item(4).Insert()

# Fill the new version column; leading apostrophe keeps values as
# text (quote-prefixed) exactly like the original run's version cells.
# Order chosen so new shared strings are created in the same order
# as the target file: version, 79.0.3945.117, 88.0.4324.96, 85.0
$ws2.Range("D1").Value = "version"
$ws2.Range("D2").Value = "'79.0.3945.117"
$ws2.Range("D6").Value = "'88.0.4324.96"
$ws2.Range("D3").Value = "'85.0"

$ws2.Range("D4").Value = "'85.0"
$ws2.Range("D5").Value = "'79.0.3945.117"
$ws2.Range("D7").Value = "'88.0.4324.96"
$ws2.Range("D8").Value = "'88.0.4324.96"
$ws2.Range("D9").Value = "'88.0.4324.96"
$ws2.Range("D10").Value = "'88.0.4324.96"
$ws2.Range("D11").Value = "'79.0.3945.117"

# multiplyLoginLogoutTest row is now executed
$ws2.Range("B6").Value = "yes"

[void]$ws2.Range("D8").Select()

# Page setup tweak recorded alongside the browser-version change
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

Write-Host "done"
